$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark. In the original document it
#    sits just before the run containing "idVersion_Current".
# -----------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# -----------------------------------------------------------------------
# 2) Locate the list-item paragraph whose whole text is "txtPage" (the
#    bullet under tblPageHistory) and rename it to "txtBody". We search
#    by exact paragraph text instead of Find/Replace so we don't touch
#    the unrelated "tblPageHistory::txtPage" substring that appears in
#    an earlier sentence.
# -----------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($paraText -eq "txtPage") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $bodyRange = $target.Range
    $bodyRange.End = $bodyRange.End - 1
    $bodyRange.Text = "txtBody"

    # ---------------------------------------------------------------
    # 3) Re-create "_GoBack" collapsed right after the renamed run
    #    (i.e. at the end of the paragraph's text, before the
    #    paragraph mark). Collapsed ranges placed directly at that
    #    offset get mis-positioned, so insert a throw-away character,
    #    bookmark across it, then delete the character -- the
    #    bookmark collapses cleanly to the correct spot.
    # ---------------------------------------------------------------
    $endPoint = $d.Range($bodyRange.End, $bodyRange.End)
    $endPoint.InsertAfter("X")

    $markerRange = $d.Range($bodyRange.End, $bodyRange.End + 1)
    $d.Bookmarks.Add("_GoBack", $markerRange)

    $bmRange = $d.Bookmarks.Item("_GoBack").Range
    $bmRange.Text = ""
}
